$wb = $excel.ActiveWorkbook

# Add a new worksheet named "dataset"
$tmp = $wb.Worksheets.Add()
$tmp.Name = "dataset"

# Move it to sit right after the existing "sample" sheet
$sampleSheet0 = $wb.Worksheets.Item("sample")
$ds0 = $wb.Worksheets.Item("dataset")
$ds0.Move($null, $sampleSheet0)

# Re-fetch fresh sheet handles now that the sheet order has changed
$sampleSheet = $wb.Worksheets.Item("sample")
$ds = $wb.Worksheets.Item("dataset")

# Header row
$ds.Range("B1").Value = "Age"
$ds.Range("C1").Value = "%"

# Section title
$ds.Range("A2").Value = "Age Groups"

# Data rows
$ds.Range("A3").Value = "18-35"
$ds.Range("B3").Value = 200
$ds.Range("C3").Value = 33.33

$ds.Range("A4").Value = "36-64"
$ds.Range("B4").Value = 200
$ds.Range("C4").Value = 33.33

$ds.Range("A5").Value = "65+"
$ds.Range("B5").Value = 200
$ds.Range("C5").Value = 33.33

# Column widths for the new sheet (target stored widths: 12.7109375 and 14.28515625 chars)
$ds.Columns.Item(1).ColumnWidth = 11.833333333333334
$ds.Columns.Item(2).ColumnWidth = 13.5

# Restore/replicate the selection on the "sample" sheet
$sampleSheet.Activate()
$sampleSheet.Range("C22").Select()

# Leave "dataset" as the active/selected tab, with its own selection
$ds.Activate()
$ds.Range("C7").Select()
